# Update the "Förändrad" (Changed) date column for all existing data rows
# (rows 2-453) from 2023-09-20 (45189) to 2023-09-21 (45190).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C453").Value = 45190

# Row 453 gains an explicit row height (matches the rest of the data rows)
# now that it is no longer the last row in the sheet.
$ws.Rows.Item(453).RowHeight = 15

# --- New row 454: A 44119-2023 -------------------------------------------
$ws.Range("A454").Value = "A 44119-2023"
$ws.Range("B454").Value = 45188
$ws.Range("B454").NumberFormat = "YYYY-MM-DD"
$ws.Range("C454").Value = 45190
$ws.Range("C454").NumberFormat = "YYYY-MM-DD"
$ws.Range("D454").Value = "VÄRMLANDS LÄN"
$ws.Range("E454").Value = "FILIPSTAD"
$ws.Range("G454").Value = 19.4
$ws.Range("H454:Q454").Value = 0
$ws.Range("R454").Value = ""
$ws.Range("R454").WrapText = $true
$ws.Rows.Item(454).RowHeight = 15

# --- New row 455: A 44503-2023 -------------------------------------------
$ws.Range("A455").Value = "A 44503-2023"
$ws.Range("B455").Value = 45189
$ws.Range("B455").NumberFormat = "YYYY-MM-DD"
$ws.Range("C455").Value = 45190
$ws.Range("C455").NumberFormat = "YYYY-MM-DD"
$ws.Range("D455").Value = "VÄRMLANDS LÄN"
$ws.Range("E455").Value = "FILIPSTAD"
$ws.Range("F455").Value = "Bergvik skog väst AB"
$ws.Range("G455").Value = 5.2
$ws.Range("H455:Q455").Value = 0
$ws.Range("R455").Value = ""
$ws.Range("R455").WrapText = $true

Write-Host "Done"
